$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts existing B..F to C..G
# and carries the style/format of the old B column (none) into new B.
$ws.Range("B1").EntireColumn.Insert()

# Capture the segment names currently sitting in column A (rows 2-20)
# before we overwrite them, then move them into the new column B.
for ($r = 2; $r -le 20; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 2).Value2 = $name
    $ws.Cells.Item($r, 2).Style = "Normal"
}

# Set the new header for column B, copying the bold/centered/bordered
# header formatting used by the other header cells.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value2 = "segments"

# Replace column A (rows 2-20) with a numeric 0-based index, keeping the
# existing "header" style that was already applied to those cells.
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}
